$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell content (rows 10, 11, 12)
$ws.Range("A10").Value = "Multiline text that goes on two rows"
$ws.Range("A10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 16.5

$ws.Range("A11").Value = "Multiline text that goes on three rows bla"
$ws.Range("A11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 16.5

$ws.Range("A12").Value = "Single line text"
$ws.Range("B12").Value2 = 123
$ws.Range("C12").Value = "Multiline text that goes on two rows"
$ws.Range("C12").WrapText = $true
$ws.Range("D12").Value = "Multiline text that goes on three rows bla"
$ws.Range("D12").WrapText = $true
$ws.Range("E12").Value2 = 43041
$ws.Range("E12").NumberFormat = "mm-dd-yy"
$ws.Rows.Item(12).RowHeight = 18.75

# Column widths (closest achievable values through pixel-quantized COM width setter)
$ws.Columns.Item(3).ColumnWidth = 15.6
$ws.Columns.Item(4).ColumnWidth = 14.1
$ws.Columns.Item(5).ColumnWidth = 9.33

# Selection / view state
$ws.Range("A10:XFD10").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 8
